# "modify W10 and add W11"
# Update the salary/task tracking sheet with this week's data:
#  - fill in date, team name, member count, team member salary
#  - replace placeholder "Member N" names with actual team member names
#  - set distributed salary as a fixed value (500) instead of a SUM formula
#  - fill in this week's completed / next week's tasks
#  - adjust a couple of row heights to fit the new content
#  - leave selection on A22, matching where the editor ended up

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header info ---
$ws.Cells.Item(3, 2).Value = "2019.11.21 - 2019.11.28"   # B3 Date
$ws.Cells.Item(4, 2).Value = "MSR Voice Input"            # B4 Team Name
$ws.Cells.Item(5, 2).Value = 5                             # B5 Total Number of Team Members

# --- Team member names (replace "Member 1".."Member 5") ---
$ws.Cells.Item(8, 1).Value = "Kunaal Sikka"
$ws.Cells.Item(9, 1).Value = "Mina Huh"
$ws.Cells.Item(10, 1).Value = "Vu Nguyen"
$ws.Cells.Item(11, 1).Value = "Nicolas Carmody"
$ws.Cells.Item(12, 1).Value = "Jonas Bokstaller"

# --- Total salary distributed: fixed value instead of the SUM formula ---
$ws.Cells.Item(14, 2).Value = 500

# --- Tasks completed this week / tasks to complete next week ---
$ws.Cells.Item(19, 1).Value = "Discuss feedbacks received during the exercise"
$ws.Cells.Item(21, 1).Value = "Update the study report"
$ws.Cells.Item(20, 1).Value = "Update the Blog"
$ws.Cells.Item(20, 1).ClearFormats()

$ws.Cells.Item(19, 2).Value = "Further update User Study report"
$ws.Cells.Item(20, 2).Value = "Plan for video prototype"

$ws.Cells.Item(22, 1).Value = "Present User Study to stakeholders"

# --- Row height tweaks to fit the refreshed content ---
$ws.Rows.Item(1).RowHeight = 91
$ws.Rows.Item(18).RowHeight = 41

# --- Selection left on A22 ---
$ws.Range("A22").Select() | Out-Null
